$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) - add row 3 for the new handback
# -----------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsOverview.Range("B3").Value = "e2e\ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-30 08:22:44"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md",
    "",
    "",
    "e2e\ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
)
$wsOverview.Range("B3").Style = "HyperLink"

$tOverview = $wsOverview.ListObjects.Item(1)
$tOverview.Resize($wsOverview.Range("A1:G3"))

# -----------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) - add row 3 for the new handback
# -----------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"

# Boolean-looking literals must be forced to text (leading apostrophe),
# then restyled to Normal so no quote-prefix flag lingers on the cell.
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("F3").Style = "Normal"

$wsZhCn.Range("G3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.4c0e51c3356ce58b7266a30e431856ebe52a3b9f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 08:22:33"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsZhCn.Range("J3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.4c0e51c3356ce58b7266a30e431856ebe52a3b9f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-30 08:23:27"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Empty-string cells: also need the apostrophe trick so a (blank) shared
# string cell is actually emitted instead of no cell at all.
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("L3").Style = "Normal"

$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("M3").Style = "Normal"

$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("N3").Style = "Normal"

$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("O3").Style = "Normal"

$wsZhCn.Range("P3").Value = "'"
$wsZhCn.Range("P3").Style = "Normal"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md",
    "",
    "",
    "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
)
$wsZhCn.Range("A3").Style = "HyperLink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/072f9a7c1c2450f1ddf9017e90c0be2aeab81f39/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md",
    "",
    "",
    "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
)
$wsZhCn.Range("I3").Style = "HyperLink"

$tZhCn = $wsZhCn.ListObjects.Item(1)
$tZhCn.Resize($wsZhCn.Range("A1:P3"))

# -----------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) - add row 3 for the new handback
# -----------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"

$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("F3").Style = "Normal"

$wsDeDe.Range("G3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.4c0e51c3356ce58b7266a30e431856ebe52a3b9f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 08:22:44"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
$wsDeDe.Range("J3").Value = "ebf5ae0e-f633-45be-8e2d-22e709d01e40.4c0e51c3356ce58b7266a30e431856ebe52a3b9f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-30 08:23:46"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("L3").Style = "Normal"

$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("M3").Style = "Normal"

$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("N3").Style = "Normal"

$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("O3").Style = "Normal"

$wsDeDe.Range("P3").Value = "'"
$wsDeDe.Range("P3").Style = "Normal"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e714b08542f8096a1a81cd0d807b6dba63bd084d/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md",
    "",
    "",
    "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
)
$wsDeDe.Range("A3").Style = "HyperLink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fd33a455d4c5c0586039d736b5d3c77f15f683a2/e2e/ebf5ae0e-f633-45be-8e2d-22e709d01e40.md",
    "",
    "",
    "ebf5ae0e-f633-45be-8e2d-22e709d01e40.md"
)
$wsDeDe.Range("I3").Style = "HyperLink"

$tDeDe = $wsDeDe.ListObjects.Item(1)
$tDeDe.Resize($wsDeDe.Range("A1:P3"))
